$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Req ID column (B2:B16) with the requirement IDs
$reqIds = @(
    "ARS_001",
    "ARS_002",
    "ARS_003",
    "ARS_004",
    "ARS_005",
    "ARS_006",
    "ARS_007",
    "ARS_008",
    "ARS_009",
    "ARS_010",
    "ARS_011",
    "ARS_012",
    "ARS_013",
    "ARS_014",
    "ARS_015"
)

for ($i = 0; $i -lt $reqIds.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $reqIds[$i]
}

# B2, B4, B6, ... B16 were filled in together (e.g. via a multi-select +
# fill), so they carry over B2's border formatting (no top edge, since it
# butts up against the header's bottom edge). Re-apply that formatting to
# the other even rows now that they hold values too.
$ws.Cells.Item(2, 2).Copy()
foreach ($row in @(4, 6, 8, 10, 12, 14, 16)) {
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Row 16's requirement text moves to row 15, and row 16 gets the
# (reworded) "track real time" requirement text in its place.
$ws.Cells.Item(15, 3).Value = "The AutoRef system shall track real time relative to the match time."
$ws.Cells.Item(16, 3).Value = "The AutoRef system shall identify the rule violations within 500 ms of actual rule violation"

# Update the selection to match the saved workbook view
$ws.Range("B2:B16").Select()
